$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BD3").Value = 151

$ws.Range("M4").Value = 1.11
$ws.Range("O4").Value = 1.63

$ws.Range("M5").Value = 1.1
$ws.Range("O5").Value = 1.54
$ws.Range("P5").Value = 2.25

$ws.Range("M6").Value = 1.08
$ws.Range("O6").Value = 1.5
$ws.Range("P6").Value = 2.37

$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 5.5
$ws.Range("O7").Value = 1.58

$ws.Range("M8").Value = 1.03
$ws.Range("N8").Value = 13
$ws.Range("O8").Value = 1.22
$ws.Range("Q8").Value = 1.77
$ws.Range("R8").Value = 1.97

$ws.Range("R10").Value = 1.62
$ws.Range("R11").Value = 1.62
